# Auto-generated edit script applying cryptos.xlsx price/volume update
# (commit: 'Updated cryptos list on Sun Mar 24 03:16:36 UTC 2024 with GitHub Actions')
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.044.42'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '3.321.66'
$ws.Range('E3').Value = '  -0.30%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '554.01'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '173.47'
$ws.Range('E6').Value = '  -0.46%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.616'
$ws.Range('E7').Value = '  +0.62%  '
$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('B9').Value = 'LidoStakedEther'
$ws.Range('C9').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D9').Value = '3.307.51'
$ws.Range('E9').Value = '  -0.55%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.169'
$ws.Range('E10').Value = '  +5.73%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.629'
$ws.Range('E11').Value = '  +1.66%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.51'
$ws.Range('E12').Value = '  -0.60%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000277'
$ws.Range('E13').Value = '  +2.78%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.02'
$ws.Range('E14').Value = '  +0.30%  '
$ws.Range('D15').Value = '3.869.34'
$ws.Range('E15').Value = '  +0.10%  '
$ws.Range('E16').Value = '  +2.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.06'
$ws.Range('E17').Value = '  -0.60%  '
$ws.Range('D18').Value = '3.343.28'
$ws.Range('E18').Value = '  +0.36%  '
$ws.Range('D19').Value = '65.244.64'
$ws.Range('E19').Value = '  +2.49%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.66'
$ws.Range('E20').Value = '  -1.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.982'
$ws.Range('E21').Value = '  +1.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '450.54'
$ws.Range('E22').Value = '  +5.76%  '
$ws.Range('E23').Value = '  +6.36%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.08'
$ws.Range('E24').Value = '  -0.08%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '86.62'
$ws.Range('E25').Value = '  +3.33%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '13.83'
$ws.Range('E26').Value = '  +6.83%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.86'
$ws.Range('E27').Value = '  +1.83%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.65'
$ws.Range('E28').Value = '  +0.73%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.56'
$ws.Range('E29').Value = '  -1.37%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '30.82'
$ws.Range('E30').Value = '  +4.47%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.55'
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.40'
$ws.Range('E32').Value = '  +0.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '566.05'
$ws.Range('E33').Value = '  -3.80%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '60.35'
$ws.Range('E34').Value = '  +3.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.106'
$ws.Range('E35').Value = '  -0.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.58'
$ws.Range('E37').Value = '  +3.37%  '
$ws.Range('E38').Value = '  -1.87%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.11'
$ws.Range('E39').Value = '  -0.43%  '
$ws.Range('D40').Value = '0.0₃0732'
$ws.Range('E40').Value = '  -1.72%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.366'
$ws.Range('E41').Value = '  +0.43%  '
$ws.Range('D42').Value = '3.057.90'
$ws.Range('E42').Value = '  -1.29%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.77'
$ws.Range('E43').Value = '  -1.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0413'
$ws.Range('E44').Value = '  +2.41%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.20'
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('E46').Value = '  +2.93%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.42'
$ws.Range('E47').Value = '  -0.11%  '
$ws.Range('E48').Value = '  +0.15%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '141.87'
$ws.Range('E49').Value = '  +6.64%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.51'
$ws.Range('E50').Value = '  -3.14%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.12'
$ws.Range('E51').Value = '  -0.32%  '
